# Update the "Checking" balance (B2) from 2000 to 300.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 300
